# Update team transition-probability matrix values (Lafayette_A)
# Changes reflect games pulled March 7 - recalculated state transition probabilities
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2081632653061224
$ws.Range("C2").Value = 0.5387755102040817
$ws.Range("J2").Value = 0.02040816326530612
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.08979591836734693
$ws.Range("B3").Value = 0.007518796992481203
$ws.Range("C3").Value = 0.03007518796992481
$ws.Range("J3").Value = 0.03759398496240601
$ws.Range("P3").Value = 0.6766917293233082
$ws.Range("S3").Value = 0.2481203007518797
$ws.Range("B6").Value = 0.039568345323741
$ws.Range("D6").Value = 0.003597122302158274
$ws.Range("F6").Value = 0.07194244604316546
$ws.Range("J6").Value = 0.2410071942446043
$ws.Range("O6").Value = 0.01079136690647482
$ws.Range("Q6").Value = 0.1510791366906475
$ws.Range("R6").Value = 0.05755395683453238
$ws.Range("S6").Value = 0.4244604316546763
$ws.Range("B7").Value = 0.08837209302325581
$ws.Range("D7").Value = 0.01395348837209302
$ws.Range("E7").Value = 0.004651162790697674
$ws.Range("F7").Value = 0.05116279069767442
$ws.Range("J7").Value = 0.1395348837209302
$ws.Range("O7").Value = 0.01395348837209302
$ws.Range("Q7").Value = 0.186046511627907
$ws.Range("R7").Value = 0.06046511627906977
$ws.Range("S7").Value = 0.4418604651162791
$ws.Range("B8").Value = 0.08314606741573034
$ws.Range("D8").Value = 0.01573033707865169
$ws.Range("E8").Value = 0.002247191011235955
$ws.Range("F8").Value = 0.06067415730337079
$ws.Range("J8").Value = 0.1303370786516854
$ws.Range("O8").Value = 0.01348314606741573
$ws.Range("Q8").Value = 0.1865168539325843
$ws.Range("R8").Value = 0.09438202247191012
$ws.Range("S8").Value = 0.4134831460674157
$ws.Range("B9").Value = 0.05164319248826291
$ws.Range("D9").Value = 0.01408450704225352
$ws.Range("F9").Value = 0.05633802816901409
$ws.Range("J9").Value = 0.107981220657277
$ws.Range("O9").Value = 0.02347417840375587
$ws.Range("Q9").Value = 0.1737089201877934
$ws.Range("R9").Value = 0.09859154929577464
$ws.Range("S9").Value = 0.4741784037558686
$ws.Range("B10").Value = 0.08400884303610906
$ws.Range("D10").Value = 0.01621223286661754
$ws.Range("F10").Value = 0.0847457627118644
$ws.Range("J10").Value = 0.09653647752394989
$ws.Range("O10").Value = 0.01915991156963891
$ws.Range("Q10").Value = 0.1834929992630803
$ws.Range("R10").Value = 0.09137803979366249
$ws.Range("S10").Value = 0.4244657332350774
$ws.Range("G11").Value = 0.1419354838709677
$ws.Range("J11").Value = 0.1032258064516129
$ws.Range("K11").Value = 0.1967741935483871
$ws.Range("L11").Value = 0.5483870967741935
$ws.Range("S11").Value = 0.00967741935483871
$ws.Range("G12").Value = 0.7790055248618785
$ws.Range("J12").Value = 0.1823204419889503
$ws.Range("L12").Value = 0.03867403314917127
$ws.Range("F13").Value = 0.02127659574468085
$ws.Range("G13").Value = 0.6595744680851063
$ws.Range("J13").Value = 0.3191489361702128
$ws.Range("F15").Value = 0.02362204724409449
$ws.Range("H15").Value = 0.1299212598425197
$ws.Range("I15").Value = 0.04724409448818898
$ws.Range("J15").Value = 0.4133858267716535
$ws.Range("K15").Value = 0.05511811023622047
$ws.Range("M15").Value = 0.003937007874015748
$ws.Range("N15").Value = 0.003937007874015748
$ws.Range("O15").Value = 0.07874015748031496
$ws.Range("S15").Value = 0.2440944881889764
$ws.Range("F16").Value = 0.007042253521126761
$ws.Range("H16").Value = 0.147887323943662
$ws.Range("I16").Value = 0.08450704225352113
$ws.Range("J16").Value = 0.4295774647887324
$ws.Range("K16").Value = 0.1197183098591549
$ws.Range("M16").Value = 0.02112676056338028
$ws.Range("O16").Value = 0.04225352112676056
$ws.Range("S16").Value = 0.147887323943662
$ws.Range("F17").Value = 0.02888888888888889
$ws.Range("H17").Value = 0.1711111111111111
$ws.Range("I17").Value = 0.1044444444444445
$ws.Range("J17").Value = 0.3955555555555555
$ws.Range("K17").Value = 0.09333333333333334
$ws.Range("M17").Value = 0.01333333333333333
$ws.Range("N17").Value = 0.002222222222222222
$ws.Range("O17").Value = 0.05111111111111111
$ws.Range("S17").Value = 0.14
$ws.Range("F18").Value = 0.02325581395348837
$ws.Range("H18").Value = 0.1302325581395349
$ws.Range("I18").Value = 0.1395348837209302
$ws.Range("J18").Value = 0.3953488372093023
$ws.Range("K18").Value = 0.08372093023255814
$ws.Range("M18").Value = 0.04186046511627907
$ws.Range("N18").Value = 0.004651162790697674
$ws.Range("O18").Value = 0.08372093023255814
$ws.Range("S18").Value = 0.09767441860465116
$ws.Range("F19").Value = 0.02389078498293516
$ws.Range("H19").Value = 0.1979522184300341
$ws.Range("I19").Value = 0.07781569965870307
$ws.Range("J19").Value = 0.3761092150170648
$ws.Range("K19").Value = 0.1037542662116041
$ws.Range("M19").Value = 0.0204778156996587
$ws.Range("N19").Value = 0.0006825938566552901
$ws.Range("O19").Value = 0.07986348122866894
$ws.Range("S19").Value = 0.1194539249146758
